$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the non-league columns first (order does not affect shared strings,
# since these values already exist in the shared string table).
$ws.Range("A3").Value = "emulator-5554"
$ws.Range("B3").Value = "Android"
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = "com.fivemobile.thescore"
$ws.Range("E3").Value = "com.fivemobile.thescore.ui.MainActivity"

$ws.Range("A4").Value = "emulator-5554"
$ws.Range("B4").Value = "Android"
$ws.Range("C4").Value = 11
$ws.Range("D4").Value = "com.fivemobile.thescore"
$ws.Range("E4").Value = "com.fivemobile.thescore.ui.MainActivity"

# New league strings must be registered in the shared string table in the
# same order they were added to the source workbook: CFL Football first,
# then MLB Baseball. Row 4 holds CFL Football and row 3 holds MLB Baseball,
# so write F4 before F3.
$ws.Range("F4").Value = "CFL Football"
$ws.Range("F3").Value = "MLB Baseball"

# Copy the border style from row 2 down to the new rows
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F4").PasteSpecial(-4122)

# Update the selected cell to match the target workbook state
$ws.Range("H6").Select()
